$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Use ClearContents (not Clear) so the header row keeps its bold/border style (s="1").
$ws.Cells.ClearContents()

# --- Header row (row 1): preserves shared-string indices 0-19 ---
$ws.Range("A1").Value = "Sending cluster"
$ws.Range("B1").Value = "Ligand symbol"
$ws.Range("C1").Value = "Receptor symbol"
$ws.Range("D1").Value = "Target cluster"
$ws.Range("E1").Value = "Ligand-expressing cells"
$ws.Range("F1").Value = "Ligand detection rate"
$ws.Range("G1").Value = "Ligand average expression value"
$ws.Range("H1").Value = "Ligand total expression value"
$ws.Range("I1").Value = "Ligand derived specificity of average expression value"
$ws.Range("J1").Value = "Ligand derived specificity of total expression value"
$ws.Range("K1").Value = "Receptor-expressing cells"
$ws.Range("L1").Value = "Receptor detection rate"
$ws.Range("M1").Value = "Receptor average expression value"
$ws.Range("N1").Value = "Receptor total expression value"
$ws.Range("O1").Value = "Receptor derived specificity of average expression value"
$ws.Range("P1").Value = "Receptor derived specificity of total expression value"
$ws.Range("Q1").Value = "Edge average expression weight"
$ws.Range("R1").Value = "Edge total expression weight"
$ws.Range("S1").Value = "Edge average expression derived specificity"
$ws.Range("T1").Value = "Edge total expression derived specificity"

# --- Data rows 2-19 ---
# row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna5"
$ws.Range("C2").Value = "Epha1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.114918
$ws.Range("H2").Value = 0.344754
$ws.Range("I2").Value = 0.04640425382421802
$ws.Range("J2").Value = 0.04640425382421801
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.086228
$ws.Range("N2").Value = 9.258683999999999
$ws.Range("O2").Value = 0.1108770174431943
$ws.Range("P2").Value = 0.1108770174431943
$ws.Range("Q2").Value = 0.354663149304
$ws.Range("R2").Value = 3.191968343736
$ws.Range("S2").Value = 0.005145165260706239
$ws.Range("T2").Value = 0.005145165260706237

# row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna5"
$ws.Range("C3").Value = "Epha1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.114918
$ws.Range("H3").Value = 0.344754
$ws.Range("I3").Value = 0.04640425382421802
$ws.Range("J3").Value = 0.04640425382421801
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.619945333333334
$ws.Range("N3").Value = 16.859836
$ws.Range("O3").Value = 0.2019043235800461
$ws.Range("P3").Value = 0.2019043235800461
$ws.Range("Q3").Value = 0.6458328778160001
$ws.Range("R3").Value = 5.812495900344
$ws.Range("S3").Value = 0.009369219479615509
$ws.Range("T3").Value = 0.009369219479615507

# row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna5"
$ws.Range("C4").Value = "Epha1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.114918
$ws.Range("H4").Value = 0.344754
$ws.Range("I4").Value = 0.04640425382421802
$ws.Range("J4").Value = 0.04640425382421801
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.355234666666667
$ws.Range("N4").Value = 10.065704
$ws.Range("O4").Value = 0.1205414547019891
$ws.Range("P4").Value = 0.1205414547019891
$ws.Range("Q4").Value = 0.385576857424
$ws.Range("R4").Value = 3.470191716816
$ws.Range("S4").Value = 0.005593636260331579
$ws.Range("T4").Value = 0.005593636260331578

# row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efna5"
$ws.Range("C5").Value = "Epha1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.114918
$ws.Range("H5").Value = 0.344754
$ws.Range("I5").Value = 0.04640425382421802
$ws.Range("J5").Value = 0.04640425382421801
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.012758666666667
$ws.Range("N5").Value = 3.038276
$ws.Range("O5").Value = 0.03638475846559173
$ws.Range("P5").Value = 0.03638475846559173
$ws.Range("Q5").Value = 0.116384200456
$ws.Range("R5").Value = 1.047457804104
$ws.Range("S5").Value = 0.001688407567170184
$ws.Range("T5").Value = 0.001688407567170184

# row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Efna5"
$ws.Range("C6").Value = "Epha1"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.114918
$ws.Range("H6").Value = 0.344754
$ws.Range("I6").Value = 0.04640425382421802
$ws.Range("J6").Value = 0.04640425382421801
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 12.86621566666667
$ws.Range("N6").Value = 38.598647
$ws.Range("O6").Value = 0.4622366263610143
$ws.Range("P6").Value = 0.4622366263610142
$ws.Range("Q6").Value = 1.478559771982
$ws.Range("R6").Value = 13.307037947838
$ws.Range("S6").Value = 0.02144974573650673
$ws.Range("T6").Value = 0.02144974573650673

# row 7
$ws.Range("A7").Value = "ECs"
$ws.Range("B7").Value = "Efna5"
$ws.Range("C7").Value = "Epha1"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.114918
$ws.Range("H7").Value = 0.344754
$ws.Range("I7").Value = 0.04640425382421802
$ws.Range("J7").Value = 0.04640425382421801
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.894313
$ws.Range("N7").Value = 5.682938999999999
$ws.Range("O7").Value = 0.06805581944816448
$ws.Range("P7").Value = 0.06805581944816448
$ws.Range("Q7").Value = 0.217690661334
$ws.Range("R7").Value = 1.959215952006
$ws.Range("S7").Value = 0.003158079519887778
$ws.Range("T7").Value = 0.003158079519887777

# row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efna5"
$ws.Range("C8").Value = "Epha1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.030023666666667
$ws.Range("H8").Value = 6.090071
$ws.Range("I8").Value = 0.819730011809897
$ws.Range("J8").Value = 0.819730011809897
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.086228
$ws.Range("N8").Value = 9.258683999999999
$ws.Range("O8").Value = 0.1108770174431943
$ws.Range("P8").Value = 0.1108770174431943
$ws.Range("Q8").Value = 6.265115880729333
$ws.Range("R8").Value = 56.386042926564
$ws.Range("S8").Value = 0.09088921881815586
$ws.Range("T8").Value = 0.09088921881815584

# row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efna5"
$ws.Range("C9").Value = "Epha1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.030023666666667
$ws.Range("H9").Value = 6.090071
$ws.Range("I9").Value = 0.819730011809897
$ws.Range("J9").Value = 0.819730011809897
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.619945333333334
$ws.Range("N9").Value = 16.859836
$ws.Range("O9").Value = 0.2019043235800461
$ws.Range("P9").Value = 0.2019043235800461
$ws.Range("Q9").Value = 11.40862203203956
$ws.Range("R9").Value = 102.677598288356
$ws.Range("S9").Value = 0.1655070335527405
$ws.Range("T9").Value = 0.1655070335527405

# row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Efna5"
$ws.Range("C10").Value = "Epha1"
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.030023666666667
$ws.Range("H10").Value = 6.090071
$ws.Range("I10").Value = 0.819730011809897
$ws.Range("J10").Value = 0.819730011809897
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.355234666666667
$ws.Range("N10").Value = 10.065704
$ws.Range("O10").Value = 0.1205414547019891
$ws.Range("P10").Value = 0.1205414547019891
$ws.Range("Q10").Value = 6.811205780553777
$ws.Range("R10").Value = 61.300852024984
$ws.Range("S10").Value = 0.09881144808644367
$ws.Range("T10").Value = 0.09881144808644367

# row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Efna5"
$ws.Range("C11").Value = "Epha1"
$ws.Range("D11").Value = "MuSCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.030023666666667
$ws.Range("H11").Value = 6.090071
$ws.Range("I11").Value = 0.819730011809897
$ws.Range("J11").Value = 0.819730011809897
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.012758666666667
$ws.Range("N11").Value = 3.038276
$ws.Range("O11").Value = 0.03638475846559173
$ws.Range("P11").Value = 0.03638475846559173
$ws.Range("Q11").Value = 2.055924061955111
$ws.Range("R11").Value = 18.503316557596
$ws.Range("S11").Value = 0.02982567848669976
$ws.Range("T11").Value = 0.02982567848669976

# row 12
$ws.Range("A12").Value = "FAPs"
$ws.Range("B12").Value = "Efna5"
$ws.Range("C12").Value = "Epha1"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.030023666666667
$ws.Range("H12").Value = 6.090071
$ws.Range("I12").Value = 0.819730011809897
$ws.Range("J12").Value = 0.819730011809897
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 12.86621566666667
$ws.Range("N12").Value = 38.598647
$ws.Range("O12").Value = 0.4622366263610143
$ws.Range("P12").Value = 0.4622366263610142
$ws.Range("Q12").Value = 26.11872230377078
$ws.Range("R12").Value = 235.068500733937
$ws.Range("S12").Value = 0.3789092351858812
$ws.Range("T12").Value = 0.3789092351858812

# row 13
$ws.Range("A13").Value = "FAPs"
$ws.Range("B13").Value = "Efna5"
$ws.Range("C13").Value = "Epha1"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.030023666666667
$ws.Range("H13").Value = 6.090071
$ws.Range("I13").Value = 0.819730011809897
$ws.Range("J13").Value = 0.819730011809897
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.894313
$ws.Range("N13").Value = 5.682938999999999
$ws.Range("O13").Value = 0.06805581944816448
$ws.Range("P13").Value = 0.06805581944816448
$ws.Range("Q13").Value = 3.845500222074333
$ws.Range("R13").Value = 34.609501998669
$ws.Range("S13").Value = 0.05578739767997609
$ws.Range("T13").Value = 0.05578739767997609

# row 14
$ws.Range("A14").Value = "MuSCs"
$ws.Range("B14").Value = "Efna5"
$ws.Range("C14").Value = "Epha1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.3315123333333334
$ws.Range("H14").Value = 0.994537
$ws.Range("I14").Value = 0.133865734365885
$ws.Range("J14").Value = 0.133865734365885
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.086228
$ws.Range("N14").Value = 9.258683999999999
$ws.Range("O14").Value = 0.1108770174431943
$ws.Range("P14").Value = 0.1108770174431943
$ws.Range("Q14").Value = 1.023122645478667
$ws.Range("R14").Value = 9.208103809308
$ws.Range("S14").Value = 0.01484263336433225
$ws.Range("T14").Value = 0.01484263336433225

# row 15
$ws.Range("A15").Value = "MuSCs"
$ws.Range("B15").Value = "Efna5"
$ws.Range("C15").Value = "Epha1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.3315123333333334
$ws.Range("H15").Value = 0.994537
$ws.Range("I15").Value = 0.133865734365885
$ws.Range("J15").Value = 0.133865734365885
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 5.619945333333334
$ws.Range("N15").Value = 16.859836
$ws.Range("O15").Value = 0.2019043235800461
$ws.Range("P15").Value = 0.2019043235800461
$ws.Range("Q15").Value = 1.863081190659111
$ws.Range("R15").Value = 16.767730715932
$ws.Range("S15").Value = 0.02702807054769015
$ws.Range("T15").Value = 0.02702807054769014

# row 16
$ws.Range("A16").Value = "MuSCs"
$ws.Range("B16").Value = "Efna5"
$ws.Range("C16").Value = "Epha1"
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.3315123333333334
$ws.Range("H16").Value = 0.994537
$ws.Range("I16").Value = 0.133865734365885
$ws.Range("J16").Value = 0.133865734365885
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 3.355234666666667
$ws.Range("N16").Value = 10.065704
$ws.Range("O16").Value = 0.1205414547019891
$ws.Range("P16").Value = 0.1205414547019891
$ws.Range("Q16").Value = 1.112301673227556
$ws.Range("R16").Value = 10.010715059048
$ws.Range("S16").Value = 0.01613637035521383
$ws.Range("T16").Value = 0.01613637035521382

# row 17
$ws.Range("A17").Value = "MuSCs"
$ws.Range("B17").Value = "Efna5"
$ws.Range("C17").Value = "Epha1"
$ws.Range("D17").Value = "MuSCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.3315123333333334
$ws.Range("H17").Value = 0.994537
$ws.Range("I17").Value = 0.133865734365885
$ws.Range("J17").Value = 0.133865734365885
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.012758666666667
$ws.Range("N17").Value = 3.038276
$ws.Range("O17").Value = 0.03638475846559173
$ws.Range("P17").Value = 0.03638475846559173
$ws.Range("Q17").Value = 0.3357419886902222
$ws.Range("R17").Value = 3.021677898212
$ws.Range("S17").Value = 0.004870672411721788
$ws.Range("T17").Value = 0.004870672411721787

# row 18
$ws.Range("A18").Value = "MuSCs"
$ws.Range("B18").Value = "Efna5"
$ws.Range("C18").Value = "Epha1"
$ws.Range("D18").Value = "Neutrophils"
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.3315123333333334
$ws.Range("H18").Value = 0.994537
$ws.Range("I18").Value = 0.133865734365885
$ws.Range("J18").Value = 0.133865734365885
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 12.86621566666667
$ws.Range("N18").Value = 38.598647
$ws.Range("O18").Value = 0.4622366263610143
$ws.Range("P18").Value = 0.4622366263610142
$ws.Range("Q18").Value = 4.265309176826556
$ws.Range("R18").Value = 38.387782591439
$ws.Range("S18").Value = 0.06187764543862637
$ws.Range("T18").Value = 0.06187764543862635

# row 19
$ws.Range("A19").Value = "MuSCs"
$ws.Range("B19").Value = "Efna5"
$ws.Range("C19").Value = "Epha1"
$ws.Range("D19").Value = "Resolving-Mac"
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.3315123333333334
$ws.Range("H19").Value = 0.994537
$ws.Range("I19").Value = 0.133865734365885
$ws.Range("J19").Value = 0.133865734365885
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 1.894313
$ws.Range("N19").Value = 5.682938999999999
$ws.Range("O19").Value = 0.06805581944816448
$ws.Range("P19").Value = 0.06805581944816448
$ws.Range("Q19").Value = 0.6279881226936666
$ws.Range("R19").Value = 5.651893104242999
$ws.Range("S19").Value = 0.009110342248300617
$ws.Range("T19").Value = 0.009110342248300615
